$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 143, shifting existing rows 143:195 down to 144:196
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new data record
$ws.Cells.Item(143, 1).Value = 10
$ws.Cells.Item(143, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(143, 3).Value = "La Araucanía"
$ws.Cells.Item(143, 4).Value = 44468
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(143, 5).Value = 9
$ws.Cells.Item(143, 6).Value = 100112044
$ws.Cells.Item(143, 7).Value = "Perejil"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 20
$ws.Cells.Item(143, 11).Value = 3300
$ws.Cells.Item(143, 12).Value = 3300
$ws.Cells.Item(143, 13).Value = 3300
$ws.Cells.Item(143, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(143, 16).Value = 1100
$ws.Cells.Item(143, 17).Value = 3
$ws.Cells.Item(143, 18).Value = "Hortaliza"
